$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-02-22 05:13:22"
$wsZhCn.Range("G5").Value = "2016-02-22 05:14:32"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-02-22 05:13:37"
$wsDeDe.Range("G5").Value = "2016-02-22 05:14:58"
